$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-looking numeric price strings to remain literal text
# (matches source data which stores these as inline strings, some
# with thousands-separator dots or trailing zeros that must survive
# verbatim instead of being normalised as Excel numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated coin snapshot values
$ws.Range("D2").Value = "28.569.69"
$ws.Range("D3").Value = "1.849.07"
$ws.Range("E3").Value = "  -4.09%  "
$ws.Range("E4").Value = "  -1.12%  "
$ws.Range("D5").Value = "336.25"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "0.4657"
$ws.Range("E7").Value = "  -3.48%  "
$ws.Range("D8").Value = "0.3900"
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07899"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "0.9766"
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "22.24"
$ws.Range("E11").Value = "  -6.28%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.839.88"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").Value = "5.812"
$ws.Range("E13").Value = "  -4.24%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "6.962"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "0.06907"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "87.80"
$ws.Range("E16").Value = "  -4.00%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001002"
$ws.Range("E18").Value = "  -3.56%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "17.02"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "28.599.37"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.388"
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "11.22"
$ws.Range("E23").Value = "  -6.08%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.155"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.104.01"
$ws.Range("E25").Value = "  -2.34%  "
$ws.Range("D26").Value = "153.49"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("D27").Value = "19.37"
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("D28").Value = "6.052"
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").Value = "2.006"
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "117.66"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "0.9656"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "0.09343"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "5.366"
$ws.Range("E33").Value = "  -4.44%  "
$ws.Range("D34").Value = "3.466"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("D35").Value = "1.347"
$ws.Range("E35").Value = "  -3.16%  "
$ws.Range("D36").Value = "0.06103"
$ws.Range("E36").Value = "  -6.69%  "
$ws.Range("D37").Value = "0.02196"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "1.162"
$ws.Range("E38").Value = "  -4.57%  "
$ws.Range("D39").Value = "0.5691"
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "7.657"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("D41").Value = "10.11"
$ws.Range("E41").Value = "  -5.83%  "
$ws.Range("D42").Value = "0.1792"
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("D44").Value = "1.253"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5366"
$ws.Range("E45").Value = "  -3.32%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "11.74"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("D47").Value = "0.07085"
$ws.Range("E47").Value = "  -6.13%  "
$ws.Range("E48").Value = "  -3.22%  "
$ws.Range("D49").Value = "112.93"
$ws.Range("E49").Value = "  -4.45%  "
$ws.Range("D50").Value = "2.344"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("E51").Value = "  -1.07%  "
